$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level updates: cryptocurrency data refreshed by scraping job
$updates = @(
    @{Row=2; D='61.985.37'; E='  -1.77%  '}
    @{Row=3; D='3.422.95'; E='  -1.02%  '}
    @{Row=4; E='  -0.02%  '}
    @{Row=5; D='578.82'; E='  -0.24%  '}
    @{Row=6; D='153.39'; E='  +3.93%  '}
    @{Row=9; D='8.05'; E='  +3.01%  '}
    @{Row=10; E='  +0.46%  '}
    @{Row=11; E='  +3.50%  '}
    @{Row=12; D='4.010.99'; E='  -1.03%  '}
    @{Row=13; E='  +0.72%  '}
    @{Row=14; D='28.52'}
    @{Row=15; D='3.441.14'; E='  -0.90%  '}
    @{Row=16; E='  -0.01%  '}
    @{Row=17; D='62.014.05'; E='  -1.72%  '}
    @{Row=18; D='6.55'}
    @{Row=19; D='14.51'; E='  +0.31%  '}
    @{Row=20; D='8.95'; E='  -3.64%  '}
    @{Row=21; D='382.97'; E='  -1.31%  '}
    @{Row=23; D='75.38'; E='  +1.14%  '}
    @{Row=25; D='3.559.39'; E='  -1.35%  '}
    @{Row=26; E='  -2.19%  '}
    @{Row=27; E='  -1.97%  '}
    @{Row=28; D='7.65'; E='  +0.35%  '}
    @{Row=29; E='  +0.64%  '}
    @{Row=30; D='2.13'; E='  -0.44%  '}
    @{Row=31; E='  -3.48%  '}
    @{Row=32; D='1.00'; E='  +0.01%  '}
    @{Row=33; D='23.26'; E='  -0.57%  '}
    @{Row=34; E='  -0.66%  '}
    @{Row=35; D='5.55'; E='  +4.25%  '}
    @{Row=36; D='1.61'; E='  +0.50%  '}
    @{Row=37; E='  -2.40%  '}
    @{Row=38; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='168.44'; E='  +0.10%  '}
    @{Row=39; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='31.06'; E='  -2.47%  '}
    @{Row=40; D='3.458.43'; E='  -1.05%  '}
    @{Row=41; D='0.0787'; E='  +2.26%  '}
    @{Row=42; D='42.68'; E='  +0.69%  '}
    @{Row=44; D='4.46'; E='  +2.09%  '}
    @{Row=45; E='  -2.82%  '}
    @{Row=47; D='2.554.08'; E='  -1.52%  '}
    @{Row=48; E='  +0.18%  '}
    @{Row=49; D='22.70'; E='  -1.19%  '}
    @{Row=50; E='  +0.02%  '}
    @{Row=51; E='  -6.06%  '}
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("B")) { $ws.Cells.Item($r, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($r, 3).Value = $u.C }
    if ($u.ContainsKey("D")) {
        # Force text storage so numeric-looking strings (e.g. "1.00", European
        # thousand-separated prices) are not coerced into floating point numbers
        $cell = $ws.Cells.Item($r, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) { $ws.Cells.Item($r, 5).Value = $u.E }
}
